# Update countries & provincias Spain
# Applies the daily-data refresh to the "Pais" sheet:
#  - Updates the "last updated" timestamp in A1
#  - Updates case statistics for several countries (some of which changed
#    relative rank and therefore swapped table rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp -------------------------------------
$ws.Range("A1").Value2 = "Datos actualizados a 31 de Mayo de 2020 a las 16:05"

# --- Helper to write a full data row (country + 7 numeric columns) -------
function Set-CountryRow($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value2 = $country
    $ws.Cells.Item($row, 2).Value2 = $b
    $ws.Cells.Item($row, 3).Value2 = $c
    $ws.Cells.Item($row, 4).Value2 = $d
    $ws.Cells.Item($row, 5).Value2 = $e
    $ws.Cells.Item($row, 6).Value2 = $f
    $ws.Cells.Item($row, 7).Value2 = $g
    $ws.Cells.Item($row, 8).Value2 = $h
}

# --- Rows whose country order swapped (rank changed) ----------------------
# India overtakes Alemania
Set-CountryRow 11 "India"    185884 4057 88546  92072 0 81 5266
Set-CountryRow 12 "Alemania" 183332 38   165200 9530  0 2  8602

# Kenia overtakes Somalia
Set-CountryRow 92 "Kenia"   1962 74 464 1435 0 0 63
Set-CountryRow 93 "Somalia" 1916 0  327 1516 0 0 73

# Belice overtakes Santa Lucia
Set-CountryRow 200 "Belice"       18 0 16 0 0 0 2
Set-CountryRow 201 "Santa Lucia"  18 0 18 0 0 0 0

# Islas Virgenes Britanicas overtakes Papua Nueva Guinea
Set-CountryRow 213 "Islas Virgenes Britanicas" 8 0 7 0 0 0 1
Set-CountryRow 214 "Papua Nueva Guinea"        8 0 8 0 0 0 0

# --- Rows with simple numeric refreshes (no reordering) -------------------
Set-CountryRow 4   "Estados Unidos" 1819788 2968 535371 1178783 0 77 105634
Set-CountryRow 52  "Oman"           11437   1014 2396   8994    0 5  47
Set-CountryRow 61  "Noruega"        8440    3    7727   477     0 0  236
Set-CountryRow 70  "Azerbaiyan"     5494    248  3428   2003    0 2  63
Set-CountryRow 101 "Sri Lanka"      1630    17   801    819     0 0  10
Set-CountryRow 110 "Libano"         1220    29   712    481     0 1  27
